# TradingModel - 2021/11/16 data update
# Adds the new TotalCapital observation for 2021-11-16 (serial date 44516)
# and keeps the "date-only" number format on the last row, shifting the
# previous last row (2021-11-15 / row 9) back to the regular datetime
# number format used by the rest of the data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember the number format currently used by the last data row (row 9),
# which is the special "date only" format reserved for the latest entry.
$lastRowDateFormat = $ws.Range("A9").NumberFormat

# Row 9 is no longer the last row, so give it the regular datetime number
# format shared by the other historical rows (e.g. row 2).
$ws.Range("A9").NumberFormat = $ws.Range("A2").NumberFormat

# Append the new data row (2021-11-16).
$ws.Range("A10").Value = 44516
$ws.Range("B10").Value = 71823.60000000001

# The new last row takes over the "date only" number format.
$ws.Range("A10").NumberFormat = $lastRowDateFormat
